$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 9 (last filled row) into rows 10 and 11 so that
# column A picks up the "interior" thin-left-border style used by the rest
# of the filled rows instead of the "first empty row" style.
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A9:F9").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 10
$ws.Range("A10").Value = "Gabriel Pereira"
$ws.Range("B10").Value = 43929
$ws.Range("C10").Value = 0.9375
$ws.Range("D10").Value = "Bataille Navale"
$ws.Range("E10").Value = "Version 1.0 fini"
$ws.Range("F10").Value = "Création de la version 1.0 terminé"

# Row 11
$ws.Range("A11").Value = "Gabriel Pereira"
$ws.Range("B11").Value = 43930
$ws.Range("C11").Value = 0.98611111111111116
$ws.Range("D11").Value = "Bataille Navale"
$ws.Range("E11").Value = "Cahier de projet fini"
$ws.Range("F11").Value = "Cahier de projet entierement remplis"

# Update selection to G11
$ws.Range("G11").Select()
